$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Match the existing table formatting (all data cells share one style:
# horizontal=left, vertical=top alignment) for the newly appended rows by
# copying the format of the last existing data row down, the same way the
# author extended the table in Excel, rather than building a brand-new
# style entry for each cell.
$ws.Range("A9:M9").Copy()
$ws.Range("A10:M13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 10 - Id 10008 "植物攻击" (plant attack)
$ws.Range("A10").Value = 10008
$ws.Range("B10").Value = "植物攻击"
$ws.Range("C10").Value = "植物的攻击"
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 2
$ws.Range("G10").Value = 0
$ws.Range("H10").Value = 1
$ws.Range("I10").Value = "hitlightning"
$ws.Range("J10").Value = "attack"
$ws.Range("K10").Value = 1667
$ws.Range("L10").Value = 820
$ws.Range("M10").Value = "hit02"

# Row 11 - Id 10009 "建筑物攻击" (building attack)
$ws.Range("A11").Value = 10009
$ws.Range("B11").Value = "建筑物攻击"
$ws.Range("C11").Value = "建筑物的攻击"
$ws.Range("D11").Value = 5
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 4
$ws.Range("G11").Value = 0
$ws.Range("H11").Value = 1
$ws.Range("I11").Value = "hitlightning"
$ws.Range("J11").Value = "attack"
$ws.Range("K11").Value = 1667
$ws.Range("L11").Value = 820
$ws.Range("M11").Value = "hit02"

# Row 12 - Id 10010 "补给物攻击" (supply attack)
$ws.Range("A12").Value = 10010
$ws.Range("B12").Value = "补给物攻击"
$ws.Range("C12").Value = "补给物的攻击"
$ws.Range("D12").Value = 1
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 2
$ws.Range("G12").Value = 0
$ws.Range("H12").Value = 1
$ws.Range("I12").Value = "hitlightning"
$ws.Range("J12").Value = "attack"
$ws.Range("K12").Value = 1667
$ws.Range("L12").Value = 820
$ws.Range("M12").Value = "hit02"

# Row 13 - Id 10011 "陷阱攻击" (trap attack)
$ws.Range("A13").Value = 10011
$ws.Range("B13").Value = "陷阱攻击"
$ws.Range("C13").Value = "陷阱的攻击"
$ws.Range("D13").Value = 2
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 2
$ws.Range("G13").Value = 0
$ws.Range("H13").Value = 1
$ws.Range("I13").Value = "hitlightning"
$ws.Range("J13").Value = "attack"
$ws.Range("K13").Value = 1667
$ws.Range("L13").Value = 820
$ws.Range("M13").Value = "hit02"

# Restore the selection to match the author's final cursor position after
# entering the new terrain/event rows.
$ws.Range("D17").Select() | Out-Null
